$d = $word.ActiveDocument

# --- Change 1: "gunicorn worktracker.wsgi:..." -> "gunicorn quantumlab.wsgi:...",
#     with the document's "_GoBack" bookmark relocated to sit right after "quantumlab".
$rng = $d.Content
$rng.Find.Execute("worktracker", $false, $false, $false, $false, $false, $true, 1, $false, "quantumlab", 2) | Out-Null

# Toggling a character-format property (net no-op) forces Word to carve this
# replaced span out into its own run instead of silently re-merging with its
# neighbours, matching the three-run split seen after the bookmark is dropped
# in the middle of the sentence.
$rng.Font.Bold = $true
$rng.Font.Bold = $false

# Word keeps a single, singleton "_GoBack" bookmark marking the last edit; Adding
# it again simply moves it here (removing it from its old location automatically).
$bmRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Change 2: collapse the three runs that spell out
#     git commit -m "   +   add something new   +   "
#     into a single run/run of text.
$rng2 = $d.Content
$rng2.Find.Execute("git commit -m ""add something new""", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Text = ""
$rng2.InsertAfter("git commit -m ""add something new""")

Write-Host "Applied gunicorn/quantumlab + _GoBack relocation + git-commit run merge."
